$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns get populated for both language
# sheets, the Status text flips from "In Translation" to
# "Handed back: in sync with en-US", and a few columns get widened so the
# new long file names are readable.
# ---------------------------------------------------------------------------

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56b3ad0006216ba98057b6400447416f0c2a4d3a/e2e/3a8b6fa8-08e8-47ef-b54a-c348951b27d6.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56b3ad0006216ba98057b6400447416f0c2a4d3a/e2e/f8081c59-8fea-459f-90e7-69a7febf5018.md"

$newStatus = "Handed back: in sync with en-US"

# Wide enough to avoid Excel's 1/6-character quantization nudging the
# rendered width away from the ~30-character target.
$wideWidth = 29.1666666666667
$fullWidth = 39.1666666666667

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlMd1, "", "", "3a8b6fa8-08e8-47ef-b54a-c348951b27d6.md")
$wsZh.Range("J2").Value = "3a8b6fa8-08e8-47ef-b54a-c348951b27d6.79a722788e714b811263aaa31499a2874c3f21e7.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-01 08:31:18"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlMd2, "", "", "f8081c59-8fea-459f-90e7-69a7febf5018.md")
$wsZh.Range("J3").Value = "f8081c59-8fea-459f-90e7-69a7febf5018.57477ad7c4c00a08a992d276840d8efd2fe712f5.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-01 08:31:18"

$wsZh.Columns.Item(3).ColumnWidth = $wideWidth
$wsZh.Columns.Item(9).ColumnWidth = $fullWidth
$wsZh.Columns.Item(10).ColumnWidth = $fullWidth

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlMd1, "", "", "3a8b6fa8-08e8-47ef-b54a-c348951b27d6.md")
$wsDe.Range("J2").Value = "3a8b6fa8-08e8-47ef-b54a-c348951b27d6.79a722788e714b811263aaa31499a2874c3f21e7.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-01 08:31:25"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlMd2, "", "", "f8081c59-8fea-459f-90e7-69a7febf5018.md")
$wsDe.Range("J3").Value = "f8081c59-8fea-459f-90e7-69a7febf5018.57477ad7c4c00a08a992d276840d8efd2fe712f5.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-01 08:31:25"

$wsDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDe.Columns.Item(10).ColumnWidth = $fullWidth
